{"js": "// Replace the 25 \"A\u00d7B=C\" answer strings in the practice-sheet table with\n// their updated values, matching the author's diff. Each old string is\n// unique in the document, so we snapshot every cell's current text first\n// (before any mutation) and then perform the substitutions by exact\n// value match. Editing via a Range.search() + insertText(..., \"Replace\")\n// on the matched range (rather than rewriting the whole cell body) keeps\n// the existing run/paragraph formatting (fonts, size, justification)\n// untouched, exactly like the diff shows.\n\nconst replacements = {\n  \"728\u00d74=2912\": \"202\u00d72=404\",\n  \"779\u00d79=7011\": \"930\u00d75=4650\",\n  \"405\u00d76=2430\": \"448\u00d73=1344\",\n  \"533\u00d75=2665\": \"923\u00d76=5538\",\n  \"365\u00d74=1460\": \"449\u00d79=4041\",\n  \"422\u00d74=1688\": \"509\u00d78=4072\",\n  \"305\u00d77=2135\": \"704\u00d75=3520\",\n  \"531\u00d74=2124\": \"874\u00d77=6118\",\n  \"849\u00d76=5094\": \"171\u00d79=1539\",\n  \"721\u00d77=5047\": \"819\u00d77=5733\",\n  \"837\u00d74=3348\": \"151\u00d78=1208\",\n  \"304\u00d75=1520\": \"937\u00d72=1874\",\n  \"426\u00d79=3834\": \"512\u00d73=1536\",\n  \"448\u00d73=1344\": \"790\u00d78=6320\",\n  \"896\u00d77=6272\": \"622\u00d73=1866\",\n  \"788\u00d73=2364\": \"406\u00d75=2030\",\n  \"331\u00d75=1655\": \"177\u00d74=708\",\n  \"914\u00d74=3656\": \"224\u00d77=1568\",\n  \"152\u00d78=1216\": \"123\u00d77=861\",\n  \"433\u00d72=866\": \"337\u00d75=1685\",\n  \"460\u00d78=3680\": \"537\u00d77=3759\",\n  \"135\u00d75=675\": \"236\u00d73=708\",\n  \"643\u00d74=2572\": \"264\u00d72=528\",\n  \"567\u00d78=4536\": \"999\u00d74=3996\",\n  \"807\u00d79=7263\": \"667\u00d79=6003\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Gather every cell from every table in the document (this file has one\n// table, but loop generically) and snapshot its current text up front.\n// `Table.columnCount` isn't available in this host, so use `values`\n// (already per-row arrays of cell text) to discover each row's width.\nconst cellBodies = [];\nfor (const table of tables.items) {\n  table.load(\"rowCount,values\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (let r = 0; r < table.rowCount; r++) {\n    const columnCount = table.values[r].length;\n    for (let c = 0; c < columnCount; c++) {\n      const cell = table.getCell(r, c);\n      cell.body.load(\"text\");\n      cellBodies.push(cell.body);\n    }\n  }\n}\nawait context.sync();\n\n// Snapshot the original text of every cell before mutating anything, so\n// a value produced by an earlier replacement (e.g. \"448\u00d73=1344\" being\n// written into one cell) is never mistaken for a still-pending source\n// value in a different cell.\nconst originalTexts = cellBodies.map((b) => b.text.replace(/\\r$/, \"\"));\n\nlet replacedCount = 0;\nfor (let i = 0; i < cellBodies.length; i++) {\n  const original = originalTexts[i];\n  const updated = replacements[original];\n  if (updated === undefined) continue;\n\n  const searchResults = cellBodies[i].search(original, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    searchResults.items[0].insertText(updated, Word.InsertLocation.replace);\n    replacedCount++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 \"A\u00d7B=C\" answer strings in the practice-sheet table with\n# their updated values, matching the author's diff.\n#\n# Every old string is unique across the document (and, trickily, one new\n# value -- \"448\u00d73=1344\" -- is also one of the *other* cells' old values).\n# A document-wide Find/Replace (or even a Find.Execute scoped to a single\n# cell's Range, which this host still resolves against the whole document\n# instead of stopping at the cell boundary) can therefore relocate to /\n# clobber the wrong cell once an earlier replacement has produced that\n# text elsewhere. To stay perfectly unambiguous we never search at all:\n# each table cell is visited by its fixed (row, column) position, its\n# current text is read and looked up directly, and -- if it matches one\n# of the known old values -- the replacement is written straight back\n# into that same `Cell.Range.Text`. Reading and writing the very same\n# cell in one step can never disturb any other cell, so no snapshot/array\n# bookkeeping is needed, and assigning `Range.Text` (rather than rebuilding\n# the cell body) keeps the existing run/paragraph formatting untouched.\n\n$replacements = @{\n  \"728\u00d74=2912\" = \"202\u00d72=404\"\n  \"779\u00d79=7011\" = \"930\u00d75=4650\"\n  \"405\u00d76=2430\" = \"448\u00d73=1344\"\n  \"533\u00d75=2665\" = \"923\u00d76=5538\"\n  \"365\u00d74=1460\" = \"449\u00d79=4041\"\n  \"422\u00d74=1688\" = \"509\u00d78=4072\"\n  \"305\u00d77=2135\" = \"704\u00d75=3520\"\n  \"531\u00d74=2124\" = \"874\u00d77=6118\"\n  \"849\u00d76=5094\" = \"171\u00d79=1539\"\n  \"721\u00d77=5047\" = \"819\u00d77=5733\"\n  \"837\u00d74=3348\" = \"151\u00d78=1208\"\n  \"304\u00d75=1520\" = \"937\u00d72=1874\"\n  \"426\u00d79=3834\" = \"512\u00d73=1536\"\n  \"448\u00d73=1344\" = \"790\u00d78=6320\"\n  \"896\u00d77=6272\" = \"622\u00d73=1866\"\n  \"788\u00d73=2364\" = \"406\u00d75=2030\"\n  \"331\u00d75=1655\" = \"177\u00d74=708\"\n  \"914\u00d74=3656\" = \"224\u00d77=1568\"\n  \"152\u00d78=1216\" = \"123\u00d77=861\"\n  \"433\u00d72=866\"  = \"337\u00d75=1685\"\n  \"460\u00d78=3680\" = \"537\u00d77=3759\"\n  \"135\u00d75=675\"  = \"236\u00d73=708\"\n  \"643\u00d74=2572\" = \"264\u00d72=528\"\n  \"567\u00d78=4536\" = \"999\u00d74=3996\"\n  \"807\u00d79=7263\" = \"667\u00d79=6003\"\n}\n\n$d = $word.ActiveDocument\n\nforeach ($table in $d.Tables) {\n  $rowCount = $table.Rows.Count\n  $colCount = $table.Columns.Count\n\n  for ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n      $cell = $table.Cell($r, $c)\n      $original = $cell.Range.Text.TrimEnd([char]7, [char]13)\n      if ($replacements.ContainsKey($original)) {\n        $cell.Range.Text = $replacements[$original]\n      }\n    }\n  }\n}\n"}
